# Weekly price-list update for "Hortaliza, Feria Lagunitas de Puerto Montt - Albahaca".
#
# A new observation (Fecha=44603 serial / 2022-02-11, Volumen=140) is inserted
# at row 63. Every existing record previously occupying rows 63..97 shifts
# down by one row (to rows 64..98), and the sheet's used range grows from
# A1:R97 to A1:R98.
#
# Implementation: walk rows 97 down to 63 and copy each row's 18 columns
# (A..R) into the row below it. Doing this in descending order means every
# source row is read before it gets overwritten by the row above it, so the
# net effect is a clean "insert a row, pushing the rest down" without
# needing a native InsertRow primitive.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 97; $r -ge 63; $r--) {
    for ($c = 1; $c -le 18; $c++) {
        $val = $ws.Cells.Item($r, $c).Value2
        $ws.Cells.Item($r + 1, $c).Value = $val
    }
}

# The freshly-written D98 (Fecha) cell has no number format yet because it
# was untouched before this script ran; match the date format used by every
# other cell in column D (e.g. D97).
$ws.Cells.Item(98, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 63 keeps its original Mercado/Región/Codreg/Categoría/Variedad/Calidad
# plus pricing/unit/origin/classification columns (copied onto itself by the
# loop above); only the new record's date and volume change.
$ws.Cells.Item(63, 4).Value = 44603
$ws.Cells.Item(63, 10).Value = 140
